# Update the subscription tier names in the header row to reflect the
# final token-based pricing discussed for the website.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "Whisker Whispers (1 Token)"
$ws.Range("E2").Value = "Meow Majesty (2 Tokens)"
$ws.Range("F2").Value = "Feline Finest (5 Tokens)"
